# Fall 2022 Schedule - Week 11 (column L) "day-after" result inputs.
# Fills in the actual game results for column L (week of 2022-11-08) which
# were previously all "A" (Available / not yet played), and nudges the
# last active-cell selection like Excel would on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Wookie Mistakes")

# row -> result value (W = Win, L = Loss, DNP = Did Not Play, NA = Not Available)
$updates = @{
    3  = "W"
    4  = "W"
    5  = "DNP"
    6  = "NA"
    7  = "L"
    8  = "W"
    9  = "NA"
    10 = "W"
    15 = "NA"
    16 = "DNP"
    17 = "W"
    18 = "W"
    19 = "W"
    20 = "DNP"
    21 = "W"
    22 = "L"
}

foreach ($row in $updates.Keys) {
    $ws.Range("L$row").Value = $updates[$row]
}

# Match the saved selection state (active cell moved from K25 to K26).
$ws.Range("K26").Select()
